$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.44954128440367
$ws.Range("C2").Value = 0.599502487562189
$ws.Range("D2").Value = 0.707865168539326
$ws.Range("E2").Value = 0.480769230769231
$ws.Range("F2").Value = 0.53448275862069

$ws.Range("B3").Value = 0.57679180887372
$ws.Range("C3").Value = 0.552112676056338
$ws.Range("D3").Value = 0.710144927536232
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 0.537051184110008
